$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new rows after row 22 (new rows 23 & 24), shifting everything
# below down by two. This turns the old single "45810" row into a 3-row
# block (22-24) describing the day's activities, mirroring the pattern
# already used for the other day blocks above (e.g. A9:A11, A12:A16, ...).
# ---------------------------------------------------------------------------
$ws.Rows.Item(23).Resize(2).Insert()

# Merge the date cell across the new 3-row block, like the other day blocks.
$ws.Range("A22:A24").Merge()

# Fill in the new activity rows.
$ws.Range("B22").Value2 = "30min"
$ws.Range("C22").Value2 = "Vérification planning"

$ws.Range("B23").Value2 = "2h30"
$ws.Range("C23").Value2 = "Test de communication port série et ethernet du convertisseur"

$ws.Range("B24").Value2 = "1h"
$ws.Range("C24").Value2 = "Recherche documentation afficheur pour la trame et de programmes pour socket"

# Row heights for the wrapped-text rows (matches the auto-computed heights
# Excel would have produced for the wrapped "Activité" text in column C).
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 45

# Update the selection to match the edited area.
$ws.Range("C23").Select()
